# Apply the "next day" update: dates and the two-digit multiplication
# answers in the table are refreshed to a new day's data.

$d = $word.ActiveDocument

function Replace-All($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# Heading date.
Replace-All "2025-10-21 Tuesday" "2025-10-22 Wednesday"

# Table answers that are unique within the document - safe to replace all
# occurrences at once.
Replace-All "35×25=875" "97×33=3201"
Replace-All "61×55=3355" "34×95=3230"
Replace-All "35×88=3080" "67×13=871"
Replace-All "45×35=1575" "66×31=2046"
Replace-All "54×35=1890" "72×80=5760"
Replace-All "66×37=2442" "25×37=925"
Replace-All "58×26=1508" "93×60=5580"
Replace-All "44×29=1276" "49×26=1274"
Replace-All "28×90=2520" "84×32=2688"
Replace-All "76×20=1520" "39×16=624"
Replace-All "16×74=1184" "85×78=6630"
Replace-All "19×43=817" "63×51=3213"
Replace-All "30×69=2070" "21×98=2058"
Replace-All "37×50=1850" "98×54=5292"
Replace-All "13×27=351" "50×36=1800"
Replace-All "43×50=2150" "47×12=564"
Replace-All "36×42=1512" "82×32=2624"
Replace-All "30×21=630" "30×87=2610"
Replace-All "69×84=5796" "44×49=2156"
Replace-All "11×85=935" "58×20=1160"
Replace-All "81×35=2835" "61×20=1220"
Replace-All "56×35=1960" "26×82=2132"
Replace-All "28×33=924" "66×24=1584"

# "80×93=7440" appears twice in the document (1st and 3rd data rows of the
# table) and must become two different values, so replace the occurrences
# one at a time, walking forward through the document without wrapping.
$r1 = $d.Content
$r1.Find.Execute("80×93=7440", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "47×28=1316", 1)

$r2 = $d.Content
$r2.Start = $r1.End
$r2.End = $d.Content.End
$r2.Find.Execute("80×93=7440", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "21×36=756", 1)
